{"js": "// Load paragraphs of the body so we can locate the first paragraph\n// (the one that starts with \"On the Insert tab, the galleries...\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\nconst newParagraphTexts = [\n  \"On the Insert tab, the galleries include items that are designed to coordinate with the overall look of your document. You can use these galleries to insert tables, headers, footers, lists, cover pages, and other document building blocks. When you create pictures, charts, or diagrams, they also coordinate with your current document look. You can easily change the formatting of selected text in the document text by choosing a look for the selected text from the Quick Styles gallery on the Home tab. You can also format text directly by using the other controls on the Home tab.\",\n  \"Most controls offer a choice of using the look from the current theme or using a format that you specify directly. To change the overall look of your document, choose new Theme elements on the Page Layout tab. To change the looks available in the Quick Style gallery, use the Change Current Quick Style Set command. Both the Themes gallery and the Quick Styles gallery provide reset commands so that you can always restore the look of your document to the original contained in your current template. On the Insert tab, the galleries include items that are designed to coordinate with the overall look of your document.\",\n  \"You can use these galleries to insert tables, headers, footers, lists, cover pages, and other document building blocks. When you create pictures, charts, or diagrams, they also coordinate with your current document look. You can easily change the formatting of selected text in the document text by choosing a look for the selected text from the Quick Styles gallery on the Home tab. You can also format text directly by using the other controls on the Home tab. Most controls offer a choice of using the look from the current theme or using a format that you specify directly.\",\n];\n\n// Insert the three new paragraphs right after the first paragraph, each one\n// right after the previous one so the final order matches the diff.\nlet anchor = firstParagraph;\nfor (const text of newParagraphTexts) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$newParagraphTexts = @(\n  \"On the Insert tab, the galleries include items that are designed to coordinate with the overall look of your document. You can use these galleries to insert tables, headers, footers, lists, cover pages, and other document building blocks. When you create pictures, charts, or diagrams, they also coordinate with your current document look. You can easily change the formatting of selected text in the document text by choosing a look for the selected text from the Quick Styles gallery on the Home tab. You can also format text directly by using the other controls on the Home tab.\",\n  \"Most controls offer a choice of using the look from the current theme or using a format that you specify directly. To change the overall look of your document, choose new Theme elements on the Page Layout tab. To change the looks available in the Quick Style gallery, use the Change Current Quick Style Set command. Both the Themes gallery and the Quick Styles gallery provide reset commands so that you can always restore the look of your document to the original contained in your current template. On the Insert tab, the galleries include items that are designed to coordinate with the overall look of your document.\",\n  \"You can use these galleries to insert tables, headers, footers, lists, cover pages, and other document building blocks. When you create pictures, charts, or diagrams, they also coordinate with your current document look. You can easily change the formatting of selected text in the document text by choosing a look for the selected text from the Quick Styles gallery on the Home tab. You can also format text directly by using the other controls on the Home tab. Most controls offer a choice of using the look from the current theme or using a format that you specify directly.\"\n)\n\n# Anchor on the first paragraph (\"On the Insert tab, the galleries...\") and\n# insert the three new paragraphs right after it, one by one, so the final\n# order matches the diff.\n$anchorIndex = 1\nforeach ($text in $newParagraphTexts) {\n    $anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n    $anchorParagraph.Range.InsertParagraphAfter()\n    $anchorIndex = $anchorIndex + 1\n    $newParagraph = $d.Paragraphs.Item($anchorIndex)\n    $newParagraph.Range.Text = $text\n}\n"}
